$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 142.9073533333333
$ws.Range("H2").Value = 428.72206
$ws.Range("I2").Value = 0.5576664151504187
$ws.Range("J2").Value = 0.5576664151504188
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.04517333333333
$ws.Range("N2").Value = 78.13552
$ws.Range("O2").Value = 0.9210237118384171
$ws.Range("P2").Value = 0.921023711838417
$ws.Range("Q2").Value = 3722.046788174578
$ws.Range("R2").Value = 33498.4210935712
$ws.Range("S2").Value = 0.5136239916494624
$ws.Range("T2").Value = 0.5136239916494624

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 142.9073533333333
$ws.Range("H3").Value = 428.72206
$ws.Range("I3").Value = 0.5576664151504187
$ws.Range("J3").Value = 0.5576664151504188
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.3302223333333333
$ws.Range("N3").Value = 0.990667
$ws.Range("O3").Value = 0.01167750336256582
$ws.Range("P3").Value = 0.01167750336256582
$ws.Range("Q3").Value = 47.19119966822444
$ws.Range("R3").Value = 424.72079701402
$ws.Range("S3").Value = 0.006512151438109043
$ws.Range("T3").Value = 0.006512151438109045

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 142.9073533333333
$ws.Range("H4").Value = 428.72206
$ws.Range("I4").Value = 0.5576664151504187
$ws.Range("J4").Value = 0.5576664151504188
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.903109
$ws.Range("N4").Value = 5.709327
$ws.Range("O4").Value = 0.06729878479901708
$ws.Range("P4").Value = 0.06729878479901708
$ws.Range("Q4").Value = 271.9682702948467
$ws.Range("R4").Value = 2447.71443265362
$ws.Range("S4").Value = 0.03753027206284734
$ws.Range("T4").Value = 0.03753027206284735

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 63.967809
$ws.Range("H5").Value = 191.903427
$ws.Range("I5").Value = 0.2496211559306514
$ws.Range("J5").Value = 0.2496211559306514
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 26.04517333333333
$ws.Range("N5").Value = 78.13552
$ws.Range("O5").Value = 0.9210237118384171
$ws.Range("P5").Value = 0.921023711838417
$ws.Range("Q5").Value = 1666.05267315856
$ws.Range("R5").Value = 14994.47405842704
$ws.Range("S5").Value = 0.2299070035886449
$ws.Range("T5").Value = 0.2299070035886448

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 63.967809
$ws.Range("H6").Value = 191.903427
$ws.Range("I6").Value = 0.2496211559306514
$ws.Range("J6").Value = 0.2496211559306514
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3302223333333333
$ws.Range("N6").Value = 0.990667
$ws.Range("O6").Value = 0.01167750336256582
$ws.Range("P6").Value = 0.01167750336256582
$ws.Range("Q6").Value = 21.123599146201
$ws.Range("R6").Value = 190.112392315809
$ws.Range("S6").Value = 0.002914951887747749
$ws.Range("T6").Value = 0.002914951887747749

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 63.967809
$ws.Range("H7").Value = 191.903427
$ws.Range("I7").Value = 0.2496211559306514
$ws.Range("J7").Value = 0.2496211559306514
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.903109
$ws.Range("N7").Value = 5.709327
$ws.Range("O7").Value = 0.06729878479901708
$ws.Range("P7").Value = 0.06729878479901708
$ws.Range("Q7").Value = 121.737713018181
$ws.Range("R7").Value = 1095.639417163629
$ws.Range("S7").Value = 0.01679920045425879
$ws.Range("T7").Value = 0.01679920045425879

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Thbs1"
$ws.Range("C8").Value = "Itga4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 49.38440333333333
$ws.Range("H8").Value = 148.15321
$ws.Range("I8").Value = 0.1927124289189298
$ws.Range("J8").Value = 0.1927124289189298
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 26.04517333333333
$ws.Range("N8").Value = 78.13552
$ws.Range("O8").Value = 0.9210237118384171
$ws.Range("P8").Value = 0.921023711838417
$ws.Range("Q8").Value = 1286.225344779911
$ws.Range("R8").Value = 11576.0281030192
$ws.Range("S8").Value = 0.1774927166003099
$ws.Range("T8").Value = 0.1774927166003099

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Thbs1"
$ws.Range("C9").Value = "Itga4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 49.38440333333333
$ws.Range("H9").Value = 148.15321
$ws.Range("I9").Value = 0.1927124289189298
$ws.Range("J9").Value = 0.1927124289189298
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.3302223333333333
$ws.Range("N9").Value = 0.990667
$ws.Range("O9").Value = 0.01167750336256582
$ws.Range("P9").Value = 0.01167750336256582
$ws.Range("Q9").Value = 16.30783289900778
$ws.Range("R9").Value = 146.77049609107
$ws.Range("S9").Value = 0.00225040003670903
$ws.Range("T9").Value = 0.00225040003670903

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Thbs1"
$ws.Range("C10").Value = "Itga4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 49.38440333333333
$ws.Range("H10").Value = 148.15321
$ws.Range("I10").Value = 0.1927124289189298
$ws.Range("J10").Value = 0.1927124289189298
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.903109
$ws.Range("N10").Value = 5.709327
$ws.Range("O10").Value = 0.06729878479901708
$ws.Range("P10").Value = 0.06729878479901708
$ws.Range("Q10").Value = 93.98390244329666
$ws.Range("R10").Value = 845.85512198967
$ws.Range("S10").Value = 0.01296931228191093
$ws.Range("T10").Value = 0.01296931228191093

